$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# Row 2: hard coal - change all values from 0.1 to 1 across B2:AF2
$ws.Range("B2:AF2").Value = 1

# Row 6: onshore wind - change all values from 0.7 to 0.5 across B6:AF6
$ws.Range("B6:AF6").Value = 0.5

# Row 13: lignite - change all values from 0.1 to 1 across B13:AF13
$ws.Range("B13:AF13").Value = 1

# Row 14: offshore wind - change all values from 0.3 to 0.1 across B14:AF14
$ws.Range("B14:AF14").Value = 0.1
